# Correção da matriz de confusão - atualizar valores da linha 2 (Testes)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Dados dos testes")

$ws.Range("B2").Value = 0.9615
$ws.Range("C2").Value = 0.9714
$ws.Range("D2").Value = 0.9722
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.9444
$ws.Range("G2").Value = 1
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0.0556
$ws.Range("J2").Value = 8
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 1
$ws.Range("M2").Value = 17
$ws.Range("N2").Value = 26
